$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6").Value = "[Python - Crawling] API 활용하여 과거 기상 관측 데이터 불러오기"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/Python-Crawling-API-%ED%99%9C%EC%9A%A9%ED%95%98%EC%97%AC-%EA%B3%BC%EA%B1%B0-%EA%B8%B0%EC%83%81-%EA%B4%80%EC%B8%A1-%EB%8D%B0%EC%9D%B4%ED%84%B0-%EB%B6%88%EB%9F%AC%EC%98%A4%EA%B8%B0"

$ws.Range("D9").Value = "[SIAI공지] Non-technical MBA 설립 허가, 한국어 모듈 프로그램 허가"
$ws.Range("E9").Value = "https://pdsi.pabii.com/siai-non-tech-mba-kr-module/#utm_source=rss&utm_medium=rss&utm_campaign=siai-non-tech-mba-kr-module"

$ws.Range("D32").Value = "Counterfactual Explanations"
$ws.Range("E32").Value = "https://dodonam.tistory.com/423"

$ws.Range("D46").Value = "Cluster headache (군발두통)"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/531"

$wb.Save()
